# feat: add 2022-Q3 data
#
# 1) Insert a brand-new worksheet "2022-Q3" right after "总计", containing the
#    fund-holdings detail table for the new quarter (shifts all the existing
#    quarter sheets one slot to the right, which Worksheets.Add(Before) does
#    for us automatically).
# 2) Insert a new top row into "总计" summarising the new quarter, shifting the
#    existing summary rows down and re-numbering the running index in column A.

function Set-TextCell($range, [string]$value) {
    # Force the cell to be stored as text (inline/shared string) even when the
    # value looks numeric (e.g. "090001", "12.69") instead of letting Excel's
    # default type-inference turn it into a number - then strip the
    # leftover NumberFormat so the cell keeps the workbook's default style.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Part 1: new "2022-Q3" worksheet with fund holdings
# ---------------------------------------------------------------------------

$anchorSheet = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($anchorSheet)
$q3.Name = "2022-Q3"

# Match page margins used throughout the rest of the workbook (0.75/0.75/1/1/0.5/0.5 in).
$q3.PageSetup.LeftMargin = 54
$q3.PageSetup.RightMargin = 54
$q3.PageSetup.TopMargin = 72
$q3.PageSetup.BottomMargin = 72
$q3.PageSetup.HeaderMargin = 36
$q3.PageSetup.FooterMargin = 36

# Clone the header row + the A-column "index" cell formatting from the
# existing "2022-Q2" sheet so the new sheet's styling matches exactly.
$anchorSheet.Range("B1:H1").Copy($q3.Range("B1:H1"))
$anchorSheet.Range("A2").Copy($q3.Range("A2"))
$anchorSheet.Range("A2").Copy($q3.Range("A3"))

# Row 2: 090001 / 大成价值增长混合
$q3.Range("A2").Value = 0
Set-TextCell $q3.Range("B2") "090001"
Set-TextCell $q3.Range("C2") "大成价值增长混合"
Set-TextCell $q3.Range("D2") "12.69"
Set-TextCell $q3.Range("E2") "64.34"
Set-TextCell $q3.Range("F2") "4.88"
Set-TextCell $q3.Range("G2") "0.6193"
$q3.Range("H2").Value = 3

# Row 3: 160919 / 大成产业升级股票（LOF）
$q3.Range("A3").Value = 1
Set-TextCell $q3.Range("B3") "160919"
Set-TextCell $q3.Range("C3") "大成产业升级股票（LOF）"
Set-TextCell $q3.Range("D3") "3.07"
Set-TextCell $q3.Range("E3") "84.12"
Set-TextCell $q3.Range("F3") "5.26"
Set-TextCell $q3.Range("G3") "0.1615"
$q3.Range("H3").Value = 5

# ---------------------------------------------------------------------------
# Part 2: prepend the 2022-Q3 row to the "总计" summary sheet
# ---------------------------------------------------------------------------

$total = $wb.Worksheets.Item("总计")

# Shift the existing data rows (2-7) down to rows 3-8, carrying formatting
# along with the values.
$total.Range("A2:D7").Copy($total.Range("A3:D8"))

# New row 2: 2022-Q3 summary
$total.Range("A2").Value = 0
Set-TextCell $total.Range("B2") "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.78

# Re-number the running index in column A for the rows that shifted down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
$total.Range("A8").Value = 6

Write-Output "2022-Q3 data added"
